$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 12).Value = "PASS"
}
